# Reorder / update the "COMPETENCES TECHNIQUES" skill lines.
#
# Before -> After (text only; paragraph count/formatting unchanged):
#   27: "Bases de donnees : SQL, MongoDB, Neo4j, Redis"              -> "Visualisation : web analytics, tableau"
#   28: "Autres : marketing, google analytics, ..."                  -> "MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit"
#   29: "Visualisation : web analytics, tableau"                     -> "Autres : marketing, google analytics, ..."
#   30: "ML/AI : Scikit-Learn, Keras, Tensorflow, ..."                -> unchanged
#   31: "MLOps : Git, DVC, Flask, Docker, ..."                       -> "Bases de donnees : SQL, MongoDB, Neo4j, Redis"
#
# We replace each paragraph's text by setting its own Range.Text directly
# (indexed via $d.Paragraphs.Item(n)) rather than using document-wide
# Find/Replace, since several of the new values (e.g. "Visualisation...",
# "MLOps...") are also old values of *other* paragraphs in this block - a
# global Find/Replace run in sequence could re-match text it just wrote.
# Scoping each replacement to its own paragraph's Range sidesteps that.

$d = $word.ActiveDocument

function Set-ParagraphText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    # Keep the paragraph mark out of the replaced text (exclude last char).
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $newText
}

Set-ParagraphText 27 "Visualisation : web analytics, tableau"
Set-ParagraphText 28 "MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit"
Set-ParagraphText 29 "Autres : marketing, google analytics, internes comme externes, presse, affichage, site centric, formats"
Set-ParagraphText 31 "Bases de données : SQL, MongoDB, Neo4j, Redis"
